{"js": "// Update the worksheet date line and the 25 two-digit-by-two-digit\n// multiplication problems (5 problem rows x 5 columns in the single table,\n// interleaved with blank answer rows) to the new values.\n\n// 1) Date line: first body paragraph, e.g. \"2025-01-03 Friday\" -> \"2025-01-04 Saturday\".\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldDate = \"2025-01-03 Friday\";\nconst newDate = \"2025-01-04 Saturday\";\nconst dateParagraph =\n  paragraphs.items.find((p) => p.text.trim() === oldDate) ?? paragraphs.items[0];\ndateParagraph.insertText(newDate, \"Replace\");\n\n// 2) Multiplication problems, addressed by (row, column) in the table so the\n// two duplicate \"26\u00d723=\" cells are each updated to their own distinct target\n// independent of text matching.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row index -> array of 5 new \"a\u00d7b=\" strings (column order left to right).\nconst newValuesByRow = {\n  0: [\"49\u00d736=\", \"36\u00d790=\", \"16\u00d733=\", \"37\u00d712=\", \"36\u00d781=\"],\n  4: [\"95\u00d770=\", \"11\u00d785=\", \"57\u00d780=\", \"85\u00d727=\", \"93\u00d754=\"],\n  9: [\"25\u00d798=\", \"46\u00d734=\", \"45\u00d713=\", \"57\u00d732=\", \"77\u00d727=\"],\n  14: [\"40\u00d760=\", \"22\u00d798=\", \"96\u00d740=\", \"99\u00d735=\", \"16\u00d785=\"],\n  19: [\"60\u00d749=\", \"39\u00d749=\", \"69\u00d721=\", \"70\u00d789=\", \"56\u00d780=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const values = newValuesByRow[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Date line: first paragraph in the document, e.g.\n#    \"2025-01-03 Friday\" -> \"2025-01-04 Saturday\".\n$find = $d.Paragraphs.Item(1).Range.Find\n$find.ClearFormatting()\n$find.Text = \"2025-01-03 Friday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2025-01-04 Saturday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Multiplication problems, addressed by (row, column) in the single table\n# so the two duplicate \"26\u00d723=\" cells are each updated to their own distinct\n# target independent of text matching.\n$t = $d.Tables.Item(1)\n\n$newValuesByRow = @{\n    1  = @(\"49\u00d736=\", \"36\u00d790=\", \"16\u00d733=\", \"37\u00d712=\", \"36\u00d781=\")\n    5  = @(\"95\u00d770=\", \"11\u00d785=\", \"57\u00d780=\", \"85\u00d727=\", \"93\u00d754=\")\n    10 = @(\"25\u00d798=\", \"46\u00d734=\", \"45\u00d713=\", \"57\u00d732=\", \"77\u00d727=\")\n    15 = @(\"40\u00d760=\", \"22\u00d798=\", \"96\u00d740=\", \"99\u00d735=\", \"16\u00d785=\")\n    20 = @(\"60\u00d749=\", \"39\u00d749=\", \"69\u00d721=\", \"70\u00d789=\", \"56\u00d780=\")\n}\n\nforeach ($rowIndex in $newValuesByRow.Keys) {\n    $values = $newValuesByRow[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
